$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-31 Wednesday", 2)

$d.Content.Find.Execute("575÷2=287, 1", $true, $false, $false, $false, $false, $true, 1, $false, "110÷6=18, 2", 2)
$d.Content.Find.Execute("174÷6=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "863÷8=107, 7", 2)
$d.Content.Find.Execute("841÷7=120, 1", $true, $false, $false, $false, $false, $true, 1, $false, "919÷5=183, 4", 2)
$d.Content.Find.Execute("513÷4=128, 1", $true, $false, $false, $false, $false, $true, 1, $false, "181÷4=45, 1", 2)
$d.Content.Find.Execute("310÷2=155, 0", $true, $false, $false, $false, $false, $true, 1, $false, "829÷9=92, 1", 2)

$d.Content.Find.Execute("632÷9=70, 2", $true, $false, $false, $false, $false, $true, 1, $false, "898÷5=179, 3", 2)
$d.Content.Find.Execute("640÷4=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "696÷8=87, 0", 2)
$d.Content.Find.Execute("151÷8=18, 7", $true, $false, $false, $false, $false, $true, 1, $false, "364÷9=40, 4", 2)
$d.Content.Find.Execute("100÷4=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "965÷7=137, 6", 2)
$d.Content.Find.Execute("983÷7=140, 3", $true, $false, $false, $false, $false, $true, 1, $false, "901÷5=180, 1", 2)

$d.Content.Find.Execute("270÷5=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "702÷9=78, 0", 2)
$d.Content.Find.Execute("731÷7=104, 3", $true, $false, $false, $false, $false, $true, 1, $false, "997÷2=498, 1", 2)
$d.Content.Find.Execute("731÷2=365, 1", $true, $false, $false, $false, $false, $true, 1, $false, "364÷2=182, 0", 2)
$d.Content.Find.Execute("247÷5=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "898÷7=128, 2", 2)
$d.Content.Find.Execute("779÷2=389, 1", $true, $false, $false, $false, $false, $true, 1, $false, "364÷3=121, 1", 2)

$d.Content.Find.Execute("454÷6=75, 4", $true, $false, $false, $false, $false, $true, 1, $false, "627÷6=104, 3", 2)
$d.Content.Find.Execute("629÷9=69, 8", $true, $false, $false, $false, $false, $true, 1, $false, "618÷7=88, 2", 2)
$d.Content.Find.Execute("854÷8=106, 6", $true, $false, $false, $false, $false, $true, 1, $false, "525÷2=262, 1", 2)
$d.Content.Find.Execute("747÷4=186, 3", $true, $false, $false, $false, $false, $true, 1, $false, "978÷3=326, 0", 2)
$d.Content.Find.Execute("636÷8=79, 4", $true, $false, $false, $false, $false, $true, 1, $false, "623÷8=77, 7", 2)

$d.Content.Find.Execute("980÷6=163, 2", $true, $false, $false, $false, $false, $true, 1, $false, "451÷9=50, 1", 2)
$d.Content.Find.Execute("600÷2=300, 0", $true, $false, $false, $false, $false, $true, 1, $false, "248÷5=49, 3", 2)
$d.Content.Find.Execute("251÷5=50, 1", $true, $false, $false, $false, $false, $true, 1, $false, "103÷9=11, 4", 2)
$d.Content.Find.Execute("900÷6=150, 0", $true, $false, $false, $false, $false, $true, 1, $false, "153÷7=21, 6", 2)
$d.Content.Find.Execute("630÷4=157, 2", $true, $false, $false, $false, $false, $true, 1, $false, "805÷3=268, 1", 2)
